# Adds the "What's New" dialog strings to the Idiomas (languages) sheet.
# New row 32 content:
#   A32 -> "¿Que hay de nuevo?"  (Spanish)
#   B32 -> "What's new?"         (English)
# Both cells keep the same centered style already used by B32 / the rest
# of the "category" rows in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new language strings.
$ws.Range("A32").Value = "¿Que hay de nuevo?"
$ws.Range("B32").Value = "What's new?"

# Keep them centered, matching the existing style used across the sheet
# (and the style B32 already had before this edit).
$ws.Range("A32:B32").HorizontalAlignment = -4108

# Move the selection to B32, and scroll the view down a bit, matching
# where the author ended up after adding the new row.
$ws.Range("B32").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
